# Bump the "想去人数" (interested-count) figures in column F for the rows
# that changed between crawls. The same events are duplicated across the
# "展览" worksheet and the "全部类型" worksheet (which aggregates every
# category), just at different row offsets, so each sheet gets its own
# row -> new-value map.
#
# NOTE: this runtime's Range/Cells ".Value" property getter (no parens)
# returns a bound-property descriptor instead of the actual cell content,
# so ".Value2" is used for both reads and writes here.

$wb = $excel.ActiveWorkbook

# "展览" sheet (row -> new value for column F)
$sheet1Updates = @{
    2  = 15236
    4  = 83
    14 = 77
    16 = 64
    17 = 19
    19 = 9032
    21 = 87
    24 = 329
    25 = 5854
    26 = 1031
    28 = 35
    29 = 83
}

# "全部类型" sheet (row -> new value for column F)
$sheet4Updates = @{
    2  = 15236
    4  = 83
    15 = 77
    17 = 64
    18 = 19
    22 = 9032
    24 = 87
    27 = 329
    28 = 5854
    29 = 1031
    31 = 35
    32 = 83
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value2 = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value2 = $sheet4Updates[$row]
}
